$d = $word.ActiveDocument

# 1. Replace the "No collaboration required" text with the new collaboration note.
$d.Content.Find.Execute("No collaboration required", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Clarification from @youngmidoriya on discord.", 2)

# 2. Fill in the Date / Time / Time Spent cells for the data row with matching formatting.
function Set-CellText($row, $col, $text) {
    $table = $word.ActiveDocument.Tables.Item(1)
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.Collapse(0)        # wdCollapseEnd -> just before the end-of-cell marker
    $rng.InsertAfter($text)

    $table2 = $word.ActiveDocument.Tables.Item(1)
    $cell2 = $table2.Cell($row, $col)
    $rng2 = $cell2.Range
    $rng2.Font.Name = "Palatino Linotype"
    $rng2.Font.Size = 10
    $rng2.Font.SizeBi = 10
}

Set-CellText 2 2 "09/06/24"
Set-CellText 2 3 "3:15pm"
Set-CellText 2 4 "1 hour"
